$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 5-7 (the two duplicate Receptor-expressing-cell breakdown rows per
# Sending cluster are collapsed: TPM recompute now aggregates the Target cluster
# receptor stats across ECs/FAPs/MuSCs into a single row per Sending cluster).
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2 updates
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002551333333333333
$ws.Range("H2").Value = 0.007654
$ws.Range("I2").Value = 0.007597293799083639
$ws.Range("J2").Value = 0.007597293799083639
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6644113333333334
$ws.Range("N2").Value = 1.993234
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.001695134781777778
$ws.Range("R2").Value = 0.015256213036
$ws.Range("S2").Value = 0.007597293799083639
$ws.Range("T2").Value = 0.007597293799083639

# Row 3 updates
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 0.05552666666666667
$ws.Range("H3").Value = 0.16658
$ws.Range("I3").Value = 0.1653458585120659
$ws.Range("J3").Value = 0.1653458585120659
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6644113333333334
$ws.Range("N3").Value = 1.993234
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.03689254663555556
$ws.Range("R3").Value = 0.33203291972
$ws.Range("S3").Value = 0.1653458585120659
$ws.Range("T3").Value = 0.1653458585120659

# Row 4 updates
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.2777433333333333
$ws.Range("H4").Value = 0.83323
$ws.Range("I4").Value = 0.8270568476888503
$ws.Range("J4").Value = 0.8270568476888503
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6644113333333334
$ws.Range("N4").Value = 1.993234
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.1845358184244444
$ws.Range("R4").Value = 1.66082236582
$ws.Range("S4").Value = 0.8270568476888503
$ws.Range("T4").Value = 0.8270568476888503
